$wb = $excel.ActiveWorkbook
$wsGeoBoundaries = $wb.Worksheets.Item("geo_boundaries")
$wsStructModel = $wb.Worksheets.Item("structuralmodel_data")

$wsGeoBoundaries.Range("F2").Value = 129.7204895019531
$wsGeoBoundaries.Range("F3").Value = 129.7204895019531
$wsGeoBoundaries.Range("F4").Value = 105.2605514526367
$wsGeoBoundaries.Range("F5").Value = 105.2605514526367
$wsGeoBoundaries.Range("F6").Value = 148.1514434814453
$wsGeoBoundaries.Range("F7").Value = 148.1514434814453
$wsGeoBoundaries.Range("F8").Value = 221.2654418945312
$wsGeoBoundaries.Range("F9").Value = 139.1471862792969
$wsGeoBoundaries.Range("F10").Value = 139.1471862792969
$wsGeoBoundaries.Range("F11").Value = 116.8295288085938
$wsGeoBoundaries.Range("F12").Value = 116.8295288085938
$wsGeoBoundaries.Range("F13").Value = 206.3394927978516
$wsGeoBoundaries.Range("F14").Value = 206.3394927978516
$wsGeoBoundaries.Range("F15").Value = 201.3811798095703
$wsGeoBoundaries.Range("F16").Value = 201.3811798095703
$wsGeoBoundaries.Range("F17").Value = 240.5259704589844
$wsGeoBoundaries.Range("F18").Value = 240.5259704589844
$wsGeoBoundaries.Range("F19").Value = 231.5464019775391
$wsGeoBoundaries.Range("F20").Value = 231.5464019775391
$wsGeoBoundaries.Range("F21").Value = 257.8213195800781
$wsGeoBoundaries.Range("F22").Value = 259.987060546875
$wsGeoBoundaries.Range("F23").Value = 168.0024566650391
$wsGeoBoundaries.Range("F24").Value = 168.0024566650391
$wsGeoBoundaries.Range("F25").Value = 256.0364990234375
$wsGeoBoundaries.Range("F26").Value = 193.2759399414062
$wsGeoBoundaries.Range("F27").Value = 193.2759399414062
$wsGeoBoundaries.Range("F28").Value = 223.0767211914062
$wsGeoBoundaries.Range("F29").Value = 223.0767211914062
$wsGeoBoundaries.Range("F30").Value = 279.4418029785156
$wsGeoBoundaries.Range("F31").Value = 210.6910400390625
$wsGeoBoundaries.Range("F32").Value = 210.6910400390625
$wsGeoBoundaries.Range("F33").Value = 284.1803588867188
$wsGeoBoundaries.Range("F34").Value = 268.218994140625
$wsGeoBoundaries.Range("F35").Value = 265.193115234375
$wsGeoBoundaries.Range("F36").Value = 210.8645935058594
$wsGeoBoundaries.Range("F37").Value = 216.2009887695312
$wsGeoBoundaries.Range("F38").Value = 307.6929626464844
$wsGeoBoundaries.Range("F39").Value = 272.0730590820312
$wsGeoBoundaries.Range("F40").Value = 226.4542846679688
$wsGeoBoundaries.Range("F41").Value = 223.0819702148438
$wsGeoBoundaries.Range("F42").Value = 211.1378479003906
$wsGeoBoundaries.Range("F43").Value = 299.9042053222656
$wsGeoBoundaries.Range("F44").Value = 295.7376098632812
$wsGeoBoundaries.Range("F45").Value = 308.8829956054688
$wsGeoBoundaries.Range("F46").Value = 223.0211029052734
$wsGeoBoundaries.Range("F47").Value = 226.7249450683594
$wsGeoBoundaries.Range("F48").Value = 233.7315521240234
$wsGeoBoundaries.Range("F49").Value = 228.9816436767578
$wsGeoBoundaries.Range("F50").Value = 276.9217224121094
$wsGeoBoundaries.Range("F51").Value = 267.8984680175781
$wsGeoBoundaries.Range("F52").Value = 288.9960021972656
$wsGeoBoundaries.Range("F53").Value = 265.5480041503906
$wsGeoBoundaries.Range("F54").Value = 281.3111877441406
$wsGeoBoundaries.Range("F55").Value = 282.0736694335938
$wsGeoBoundaries.Range("F56").Value = 258.8345336914062
$wsGeoBoundaries.Range("F57").Value = 271.9231262207031
$wsGeoBoundaries.Range("F58").Value = 221.2654418945312

$wsStructModel.Range("D2").Value = 201.3811798095703
$wsStructModel.Range("D3").Value = 201.3811798095703
$wsStructModel.Range("D4").Value = 201.3811798095703
$wsStructModel.Range("D5").Value = 240.5259704589844
$wsStructModel.Range("D6").Value = 240.5259704589844
$wsStructModel.Range("D7").Value = 240.5259704589844
$wsStructModel.Range("D8").Value = 231.5464019775391
$wsStructModel.Range("D9").Value = 231.5464019775391
$wsStructModel.Range("D10").Value = 231.5464019775391
$wsStructModel.Range("D11").Value = 257.8213195800781
$wsStructModel.Range("D12").Value = 232.8213195800781
$wsStructModel.Range("D13").Value = 167.8213195800781
$wsStructModel.Range("D14").Value = 259.987060546875
$wsStructModel.Range("D15").Value = 174.987060546875
$wsStructModel.Range("D16").Value = 122.987060546875
$wsStructModel.Range("D17").Value = 168.0024566650391
$wsStructModel.Range("D18").Value = 168.0024566650391
$wsStructModel.Range("D19").Value = 168.0024566650391
$wsStructModel.Range("D20").Value = 256.0364990234375
$wsStructModel.Range("D21").Value = 40.0364990234375
$wsStructModel.Range("D22").Value = -39.9635009765625
$wsStructModel.Range("D23").Value = 193.2759399414062
$wsStructModel.Range("D24").Value = 162.2759399414062
$wsStructModel.Range("D25").Value = 193.2759399414062
$wsStructModel.Range("D26").Value = 223.0767211914062
$wsStructModel.Range("D27").Value = 156.0767211914062
$wsStructModel.Range("D28").Value = 223.0767211914062
$wsStructModel.Range("D29").Value = 279.4418029785156
$wsStructModel.Range("D30").Value = 185.4418029785156
$wsStructModel.Range("D31").Value = 130.4418029785156
$wsStructModel.Range("D32").Value = 210.6910400390625
$wsStructModel.Range("D33").Value = 210.6910400390625
$wsStructModel.Range("D34").Value = 210.6910400390625
$wsStructModel.Range("D35").Value = 284.1803588867188
$wsStructModel.Range("D36").Value = 101.1803588867188
$wsStructModel.Range("D37").Value = 43.18035888671881
$wsStructModel.Range("D38").Value = 268.218994140625
$wsStructModel.Range("D39").Value = 109.218994140625
$wsStructModel.Range("D40").Value = 39.218994140625
$wsStructModel.Range("D41").Value = 265.193115234375
$wsStructModel.Range("D42").Value = 18.193115234375
$wsStructModel.Range("D43").Value = 210.8645935058594
$wsStructModel.Range("D44").Value = 27.8645935058594
$wsStructModel.Range("D45").Value = 216.2009887695312
$wsStructModel.Range("D46").Value = 64.20098876953119
$wsStructModel.Range("D47").Value = 307.6929626464844
$wsStructModel.Range("D48").Value = 103.6929626464844
$wsStructModel.Range("D49").Value = 71.69296264648438
$wsStructModel.Range("D50").Value = 272.0730590820312
$wsStructModel.Range("D51").Value = 77.07305908203119
$wsStructModel.Range("D52").Value = 50.07305908203119
$wsStructModel.Range("D53").Value = 226.4542846679688
$wsStructModel.Range("D54").Value = 31.45428466796881
$wsStructModel.Range("D55").Value = 223.0819702148438
$wsStructModel.Range("D56").Value = 86.08197021484381
$wsStructModel.Range("D57").Value = 211.1378479003906
$wsStructModel.Range("D58").Value = 44.1378479003906
$wsStructModel.Range("D59").Value = 17.1378479003906
$wsStructModel.Range("D60").Value = 299.9042053222656
$wsStructModel.Range("D61").Value = 100.9042053222656
$wsStructModel.Range("D62").Value = 70.90420532226562
$wsStructModel.Range("D63").Value = 295.7376098632812
$wsStructModel.Range("D64").Value = -138.2623901367188
$wsStructModel.Range("D65").Value = -232.2623901367188
$wsStructModel.Range("D66").Value = 308.8829956054688
$wsStructModel.Range("D67").Value = 201.8829956054688
$wsStructModel.Range("D68").Value = 33.88299560546881
$wsStructModel.Range("D69").Value = 223.0211029052734
$wsStructModel.Range("D70").Value = 113.0211029052734
$wsStructModel.Range("D71").Value = 226.7249450683594
$wsStructModel.Range("D72").Value = 124.7249450683594
$wsStructModel.Range("D73").Value = 233.7315521240234
$wsStructModel.Range("D74").Value = 75.73155212402341
$wsStructModel.Range("D75").Value = 228.9816436767578
$wsStructModel.Range("D76").Value = 121.9816436767578
$wsStructModel.Range("D77").Value = 276.9217224121094
$wsStructModel.Range("D78").Value = -25.07827758789062
$wsStructModel.Range("D79").Value = 288.9960021972656
$wsStructModel.Range("D80").Value = -10.00399780273438
$wsStructModel.Range("D81").Value = 265.5480041503906
$wsStructModel.Range("D82").Value = -112.4519958496094
$wsStructModel.Range("D83").Value = -214.4519958496094
$wsStructModel.Range("D84").Value = 281.3111877441406
$wsStructModel.Range("D85").Value = 211.3111877441406
$wsStructModel.Range("D86").Value = 282.0736694335938
$wsStructModel.Range("D87").Value = 177.0736694335938
$wsStructModel.Range("D88").Value = 258.8345336914062
$wsStructModel.Range("D89").Value = -32.16546630859381
$wsStructModel.Range("D90").Value = -98.16546630859381
